$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price (D) and volume-change (E) columns for rows 2-51 with
# refreshed figures from the latest scrape.

# The Price column stores text (e.g. "29.189.02", "1.0000") rather than numbers -
# force every D-column cell we touch to keep a Text format so COM does not
# reinterpret digit-and-dot strings as numbers and mangle them (stripping trailing
# zeros, re-grouping thousands separators, etc).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.189.02"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.00"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.39"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6188"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07369"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2913"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.08"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.46"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.948"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6646"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.23"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008914"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.162.06"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.070.98"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.09"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.349"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.89"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.491"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.63"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05917"
$ws.Range("E30").Value = "  +5.98%  "
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.854"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7292"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.604"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.845"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219.86"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01746"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.279"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9180"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.86"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.976.51"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.78"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5087"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.119"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4017"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  +2.00%  "
